# "separate dept from affiliations"
#
# PI hours: add a new "app" column holding the full affiliation list that
# used to live in "dept", and shrink "dept" itself down to the PI's single
# primary department.
#
# "dept hours" is renamed to "department hours" and its data is rebuilt to
# aggregate hours by the new single-department values.
#
# A brand-new "unit(accumulative) hours" sheet is appended, carrying the
# original ("dept hours") accumulative-affiliation aggregation that used to
# live on the second tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "PI hours" - add the "app" column, trim "dept" to one value
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cell F1, styled like the other header cells (B1:E1)
$ws1.Range("F1").Value = "app"
$ws1.Range("E1").Copy()
$ws1.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2: Romit Roy Choudhury
$ws1.Range("F2").Value = "['ECE', 'CSL']"
$ws1.Range("E2").Value = "ECE"

# Row 3: Naira Hovakimyan
$ws1.Range("F3").Value = "['ME', 'AE', 'CSL']"
$ws1.Range("E3").Value = "ME"

# Row 4: Paul G Kwiat
$ws1.Range("F4").Value = "['PHYS', 'ECE']"
$ws1.Range("E4").Value = "PHYS"

# Row 5: Sayan Mitra
$ws1.Range("F5").Value = "['ECE', 'CSL']"
$ws1.Range("E5").Value = "ECE"

# ---------------------------------------------------------------------
# Sheet 2: rename "dept hours" -> "department hours", rebuild its data
# to aggregate by the (now single-valued) department. It now only needs
# 3 data rows (was 5), so drop the trailing two.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A5:D6").Delete()
$ws2.Name = "department hours"

$ws2.Range("B1").Value = "dept"

$ws2.Range("B2").Value = "ECE"
$ws2.Range("C2").Value = 53
$ws2.Range("D2").Value = 55.78947368421053

$ws2.Range("B3").Value = "ME"
$ws2.Range("C3").Value = 35
$ws2.Range("D3").Value = 36.8421052631579

$ws2.Range("B4").Value = "PHYS"
$ws2.Range("C4").Value = 7
$ws2.Range("D4").Value = 7.368421052631579

# ---------------------------------------------------------------------
# Sheet 3 (new): "unit(accumulative) hours" - carries the old "dept hours"
# accumulative-affiliation data (6 rows, same shape the second tab used
# to have before this edit).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "unit(accumulative) hours"

$ws3.Range("B1").Value = "unit(accumulative)"
$ws3.Range("C1").Value = "hours"
$ws3.Range("D1").Value = "percentage"

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "CSL"
$ws3.Range("C2").Value = 88
$ws3.Range("D2").Value = 39.11111111111111

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "ECE"
$ws3.Range("C3").Value = 60
$ws3.Range("D3").Value = 26.66666666666667

$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "ME"
$ws3.Range("C4").Value = 35
$ws3.Range("D4").Value = 15.55555555555556

$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "AE"
$ws3.Range("C5").Value = 35
$ws3.Range("D5").Value = 15.55555555555556

$ws3.Range("A6").Value = 4
$ws3.Range("B6").Value = "PHYS"
$ws3.Range("C6").Value = 7
$ws3.Range("D6").Value = 3.111111111111111

# Match the header/index-column styling used elsewhere in the workbook
# (bold + centered + bordered header row, bold + bordered index column).
$ws1.Range("B1:D1").Copy()
$ws3.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("A2:A5").Copy()
$ws3.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("A2").Copy()
$ws3.Range("A6").PasteSpecial(-4122)  # xlPasteFormats

# Leave the view on the originally-active sheet ("PI hours") rather than
# the freshly-added one.
$ws1.Activate()
